$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix casing of header labels: metadata4Ing_* -> metadata4ing_*
$ws.Range("B1").Value = "metadata4ing_IRI"
$ws.Range("C1").Value = "metadata4ing_DESC"

# Add new header column F1, cloning formatting (bold/centered/bordered) from E1
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "REX_DEF"

# Update the description text for the metadata4ing concept (C2)
$ws.Range("C2").Value = "{'label': None, 'prefLabel': 'Association', 'altLabel': None, 'name': 'Association'}"

# Add new data cell F2
$ws.Range("F2").Value = "[]"
